$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(36274, 52505563)
    3 = @(88126, 129291228)
    4 = @(30196, 44748218)
    5 = @(8376, 12452176)
    6 = @(1831, 2723069)
    7 = @(137, 200593)
    11 = @(39615, 53836906)
    12 = @(9318, 13482865)
    13 = @(25248, 37048192)
    14 = @(8086, 12007519)
    15 = @(2078, 3092150)
    16 = @(386, 568123)
    17 = @(30, 45000)
    19 = @(9812, 13031821)
    20 = @(13008, 18791888)
    21 = @(30894, 45371856)
    22 = @(10007, 14881965)
    23 = @(2554, 3801263)
    24 = @(466, 693345)
    25 = @(33, 48953)
    26 = @(11311, 15142280)
    27 = @(7383, 10699316)
    28 = @(21901, 32151618)
    29 = @(7596, 11305513)
    30 = @(1900, 2835544)
    31 = @(335, 499915)
    33 = @(8032, 10628942)
    34 = @(3043, 4387109)
    35 = @(7435, 10865041)
    36 = @(2998, 4441288)
    37 = @(790, 1177763)
    40 = @(2294, 3096957)
    41 = @(16663, 24108494)
    42 = @(49624, 72783191)
    43 = @(18518, 27507707)
    44 = @(5411, 8060312)
    45 = @(1101, 1642292)
    46 = @(58, 85348)
    49 = @(16120, 21501776)
    50 = @(1839, 2670175)
    51 = @(6389, 9400983)
    52 = @(2211, 3301958)
    53 = @(718, 1072305)
    54 = @(169, 250333)
    56 = @(6013, 8295737)
    57 = @(778, 1139425)
    58 = @(2000, 2967826)
    59 = @(813, 1211953)
    60 = @(271, 406258)
    61 = @(65, 97500)
    63 = @(1171, 1658558)
    64 = @(14871, 21492379)
    65 = @(43584, 63813422)
    66 = @(15324, 22782406)
    67 = @(4430, 6599060)
    68 = @(879, 1308096)
    69 = @(70, 102330)
    71 = @(14628, 19328332)
    72 = @(48580, 70734459)
    73 = @(139326, 205363069)
    74 = @(60906, 90786113)
    75 = @(19370, 28946079)
    76 = @(4425, 6611472)
    77 = @(242, 358170)
    83 = @(48133, 65743256)
    84 = @(4363, 6326271)
    85 = @(11082, 16286094)
    86 = @(3765, 5610915)
    87 = @(1308, 1954489)
    88 = @(278, 414512)
    91 = @(5063, 6823393)
    92 = @(1479, 2137783)
    93 = @(4853, 7148054)
    94 = @(1855, 2763999)
    95 = @(658, 985960)
    96 = @(168, 251113)
    99 = @(3266, 4333166)
    100 = @(564, 841464)
    101 = @(316, 472130)
    102 = @(112, 168000)
    105 = @(10424, 15141890)
    106 = @(28510, 41900426)
    107 = @(9554, 14208057)
    108 = @(2602, 3879307)
    109 = @(459, 685982)
    112 = @(9476, 12544900)
    113 = @(29377, 42399516)
    114 = @(64406, 94306401)
    115 = @(20851, 30998479)
    116 = @(5855, 8723780)
    117 = @(1064, 1590506)
    118 = @(67, 97920)
    121 = @(24916, 33333054)
    122 = @(34440, 49747000)
    123 = @(74152, 108512282)
    124 = @(23081, 34268153)
    125 = @(6155, 9148291)
    126 = @(1149, 1707911)
    130 = @(30360, 40398179)
    131 = @(12828, 18575620)
    132 = @(31518, 46315522)
    133 = @(11221, 16673059)
    134 = @(2871, 4281081)
    135 = @(460, 683990)
    138 = @(10484, 14010177)
    139 = @(33673, 48661631)
    140 = @(78592, 115199792)
    141 = @(23635, 35136925)
    142 = @(6160, 9194586)
    143 = @(1359, 2022599)
    144 = @(73, 109130)
    146 = @(28128, 38060020)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}
